# Updated PvsI with model fitting: refresh volume/area-derived rate columns
# (T, V, Z, AB, AC, AD) for rows 2-15 of the active sheet with new model
# output, and switch the per-area output unit from cm^2 to m^2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUnit = "umolO2/min/m2"

# row => @(T, V, Z, AB, AD)   AB/AD can be numeric, the string "Inf", or
# $null (meaning: clear the cell entirely, as in row 15 of the diff).
$rows = @{
    2  = @(0.1450048780487805, 0.0002448603057459146, -0.219062907004309,  -894.6444232231951, -894.6444232231951)
    3  = @(0.1492487804878049, 0.0001488973818309612, -0.2311228863950921, -1552.229351201616, -1552.229351201616)
    4  = @(0.1469268292682927, 0.0002222807942365138, -0.1717472054618187, -772.6587717653832,  -772.6587717653832)
    5  = @(0.1418926829268293, 0.0002529432437181515, -0.2428954027474859, -960.2763022132283,  -960.2763022132283)
    6  = @(0.1446439024390244, 0.0001851607801792304, -0.2258674431593748, -1219.844952806645,  -1219.844952806645)
    7  = @(0.1429268292682927, 0.0003232296608680373, -0.2340086886050746, -723.9703434908829,  -723.9703434908829)
    8  = @(0.1544,             0,                      0.001412416365969315, "Inf",             "Inf")
    9  = @(0.1450048780487805, 0.0002448603057459146,  0.2353534575927554,  961.1744005456558,   961.1744005456558)
    10 = @(0.1492487804878049, 0.0001488973818309612,  0.2937119032397754,  1972.579367266631,   1972.579367266631)
    11 = @(0.1469268292682927, 0.0002222807942365138,  0.1223396063125632,  550.3831616796813,   550.3831616796813)
    12 = @(0.1418926829268293, 0.0002529432437181515,  0.2762386004742184,  1092.097169363513,   1092.097169363513)
    13 = @(0.1446439024390244, 0.0001851607801792304,  0.2791288544710472,  1507.494482367478,   1507.494482367478)
    14 = @(0.1429268292682927, 0.0003232296608680373,  0.2429251293160247,  751.5558091533005,   751.5558091533005)
    15 = @(0.1544,             0,                      0,                   $null,               $null)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $tVal  = $vals[0]
    $vVal  = $vals[1]
    $zVal  = $vals[2]
    $abVal = $vals[3]
    $adVal = $vals[4]

    $ws.Cells.Item($r, 20).Value = $tVal   # T = volume
    $ws.Cells.Item($r, 22).Value = $vVal   # V = area
    $ws.Cells.Item($r, 26).Value = $zVal   # Z = rate.abs

    if ($null -eq $abVal) {
        $ws.Cells.Item($r, 28).ClearContents()   # AB = rate.a.spec
    } else {
        $ws.Cells.Item($r, 28).Value = $abVal
    }

    $ws.Cells.Item($r, 29).Value = $newUnit   # AC = output.unit

    if ($null -eq $adVal) {
        $ws.Cells.Item($r, 30).ClearContents()   # AD = rate.output
    } else {
        $ws.Cells.Item($r, 30).Value = $adVal
    }
}
